$wb = $excel.ActiveWorkbook

# --- "Heat Generators" sheet: zero out renewable factor (column K) for several rows ---
$wsHeat = $wb.Worksheets.Item("Heat Generators")

$rows = @(4,5,6,7,8,9,13,14,15,16,17,22,23,24,25,26,27,28,29)
foreach ($r in $rows) {
    $wsHeat.Cells.Item($r, 11).Value = 0
}

# --- "financal and other parameteres" sheet: update B3 ---
$wsFin = $wb.Worksheets.Item("financal and other parameteres")
$wsFin.Range("B3").Value = 0.05

# --- Update selections on each sheet ---
$wsHeat.Range("K3:K29").Select() | Out-Null
$wsHeat.Activate() | Out-Null

$wsFin.Range("D13").Select() | Out-Null
$wsFin.Activate() | Out-Null
